$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 284, pushing existing rows 284-378 down to 286-380.
$ws.Rows("284:285").Insert()

# Populate new row 284 (all columns A-R)
$ws.Cells.Item(284,1).Value = 9
$ws.Cells.Item(284,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(284,3).Value = "Metropolitana"
$ws.Cells.Item(284,4).Value = 44588
$ws.Cells.Item(284,5).Value = 13
$ws.Cells.Item(284,6).Value = 100112031
$ws.Cells.Item(284,7).Value = "Poroto verde"
$ws.Cells.Item(284,8).Value = "Magnum"
$ws.Cells.Item(284,9).Value = "Primera"
$ws.Cells.Item(284,10).Value = 70
$ws.Cells.Item(284,11).Value = 36000
$ws.Cells.Item(284,12).Value = 38000
$ws.Cells.Item(284,13).Value = 37000
$ws.Cells.Item(284,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(284,15).Value = "Región Metropolitana"
$ws.Cells.Item(284,16).Value = 1480
$ws.Cells.Item(284,17).Value = 25
$ws.Cells.Item(284,18).Value = "Hortaliza"

# Populate new row 285 (all columns A-R)
$ws.Cells.Item(285,1).Value = 9
$ws.Cells.Item(285,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(285,3).Value = "Metropolitana"
$ws.Cells.Item(285,4).Value = 44588
$ws.Cells.Item(285,5).Value = 13
$ws.Cells.Item(285,6).Value = 100112031
$ws.Cells.Item(285,7).Value = "Poroto verde"
$ws.Cells.Item(285,8).Value = "Sin especificar"
$ws.Cells.Item(285,9).Value = "Primera"
$ws.Cells.Item(285,10).Value = 43
$ws.Cells.Item(285,11).Value = 46000
$ws.Cells.Item(285,12).Value = 48000
$ws.Cells.Item(285,13).Value = 47023
$ws.Cells.Item(285,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(285,15).Value = "Provincia del Elquí"
$ws.Cells.Item(285,16).Value = 1881
$ws.Cells.Item(285,17).Value = 25
$ws.Cells.Item(285,18).Value = "Hortaliza"
